$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 2
$ws.Range("F4").Value = -3
$ws.Range("H4").Value = 46

$ws.Range("D4").Select()
